$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a row's "weekly" data block.
$cols = @(4, 13, 14, 15, 16, 19)   # D, M, N, O, P, S

# Rows affected by the reshuffle (row 1 is the header, row 14 is untouched).
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 15, 16, 17, 18)

# Snapshot the original values for each affected row/column before
# overwriting anything (the mapping below is a permutation, not an
# independent set of edits, so sources must be read before any writes).
$orig = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# Destination row -> source row (i.e. destination row now holds what used
# to be in the source row).
$mapping = @{
    2  = 10
    3  = 8
    4  = 17
    5  = 11
    6  = 15
    7  = 18
    8  = 5
    9  = 12
    10 = 6
    11 = 7
    12 = 2
    13 = 9
    15 = 16
    16 = 4
    17 = 13
    18 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
